$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$ws.Cells.Item(72, 1).Value = 11
$ws.Cells.Item(72, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(72, 3).Value = "Bíobío"
$ws.Cells.Item(72, 4).Value = 44778
$ws.Cells.Item(72, 5).Value = 8
$ws.Cells.Item(72, 6).Value = 100112032
$ws.Cells.Item(72, 7).Value = "Zapallo italiano"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 170
$ws.Cells.Item(72, 11).Value = 19000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 19529
$ws.Cells.Item(72, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 325
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"
